$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3 currently describes a markdown "folder example" entry:
#   A3 = folder_1_md, B3 = md, C3 = example_1, D3 = data/md/folder/example_1.md
# Rename/replace it with the new "tourisme_exemple" markdown entry.
# B3 ("md") stays the same value.
$ws.Range("A3").Value = "tourisme_exemple"
$ws.Range("C3").Value = "Tourisme exemple"
$ws.Range("D3").Value = "data/md/tourisme_exemple.md"

# Update the active selection shown in the saved view.
$ws.Range("D13").Select()
